$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = "Feat: Merge Height fix logic into main roster scrape code"
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = $ws.Range("B13").Value2
$ws.Range("C14").Value = "M"

$ws.Range("C14").Select()
